$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.526.44"
$ws.Range("E2").Value = "  +2.20%  "

# Row 3
$ws.Range("D3").Value = "3.596.34"
$ws.Range("E3").Value = "  +5.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "651.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.93%  "

# Row 7
$ws.Range("E7").Value = "  +5.54%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.405"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.30%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("E10").Value = "  +2.63%  "

# Row 11
$ws.Range("D11").Value = "3.591.71"
$ws.Range("E11").Value = "  +4.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.73"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.68%  "

# Row 13
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.32"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("D15").Value = "4.282.01"
$ws.Range("E15").Value = "  +5.24%  "

# Row 16
$ws.Range("D16").Value = "95.463.99"
$ws.Range("E16").Value = "  +2.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000254"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.98%  "

# Row 18
$ws.Range("D18").Value = "3.597.15"
$ws.Range("E18").Value = "  +5.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.40%  "

# Row 20
$ws.Range("E20").Value = "  +8.81%  "

# Row 21
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.486"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "509.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.37%  "

# Row 25
$ws.Range("E25").Value = "  +6.10%  "

# Row 26
$ws.Range("E26").Value = "  -0.26%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.74%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.24%  "

# Row 29
$ws.Range("D29").Value = "3.799.89"
$ws.Range("E29").Value = "  +5.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +12.90%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.93%  "

# Row 32
$ws.Range("E32").Value = "  +1.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("E35").Value = "  +3.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.90"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.559"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.21"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.32%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "572.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.18%  "

# Row 41
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
$ws.Range("E42").Value = "  +0.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.920"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.77"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.36%  "

# Row 45
$ws.Range("E45").Value = "  -0.30%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.52%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.30"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +34.61%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.99%  "

# Row 49
$ws.Range("E49").Value = "  +1.00%  "

# Row 50
$ws.Range("E50").Value = "  -3.71%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.74"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.20%  "
